$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cell B11 needs to hold the literal text "1" (it currently holds the text
# "R40") while its existing cell style (s="23") must stay exactly as-is.
# A plain `.Value = "1"` assignment gets reinterpreted as a *number* (or,
# if we force text via NumberFormat first, it bumps the cell onto a new
# style), so neither approach preserves the original style untouched.
#
# Instead we stage the literal text in an already-blank, unmerged cell
# that sits inside the sheet's existing used range (B5), using a formula
# that evaluates to text so no NumberFormat change (and therefore no new
# style) is ever needed. We then copy *only the value* over to B11, which
# leaves B11's style completely alone, and finally restore B5 back to its
# original empty state.
$scratch = $ws.Range("B5")
$scratch.Formula = '="1"'

$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$scratch.ClearContents()
